$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header text: D4 "HW" -> "Docs", F4 "Quiz" -> "Quiz(13)"
$ws.Range("D4").Value = "Docs"
$ws.Range("F4").Value = "Quiz(13)"

# Add grade values for quiz entries (rows 9 & 10, column C)
$ws.Range("C9").Value = 1
$ws.Range("C10").Value = 1

# Update selected cell to match final state
$ws.Range("C11").Select()
